# Auto-generated Excel COM-interop script to apply the Siren_Profits market-data refresh
# described by the commit diff. Each worksheet (by index, matching workbook.xml sheet order)
# receives direct cell value updates; cells removed by the diff are cleared via ClearContents().

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H15").Value = 1201.8182
$ws.Range("I15").Value = 1201.8182
$ws.Range("K15").Value = 3605.4546
$ws.Range("M15").Value = -3436.4546
$ws.Range("H18").Value = 3352.7144
$ws.Range("I18").Value = 3444.2964
$ws.Range("K18").Value = 3444.2964
$ws.Range("M18").Value = -3160.2964
$ws.Range("H74").Value = 8994.75
$ws.Range("J74").Value = 8994.75
$ws.Range("L74").Value = 8994.75
$ws.Range("N74").Value = -10866.75
$ws.Range("H77").Value = 8994.75
$ws.Range("J77").Value = 8994.75
$ws.Range("L77").Value = 44973.75
$ws.Range("N77").Value = -54333.75
$ws.Range("H88").Value = 12876.25
$ws.Range("J88").Value = 6051.5
$ws.Range("L88").Value = 6051.5
$ws.Range("N88").Value = -6863.5
$ws.Range("H91").Value = 12876.25
$ws.Range("J91").Value = 6051.5
$ws.Range("L91").Value = 6051.5
$ws.Range("N91").Value = -8859.5
$ws.Range("H92").Value = 1141.2
$ws.Range("I92").Value = 582.0714
$ws.Range("J92").Value = 2445.8333
$ws.Range("K92").Value = 582.0714
$ws.Range("L92").Value = 2445.8333
$ws.Range("M92").Value = 665.9286
$ws.Range("N92").Value = -4941.8333
$ws.Range("H94").Value = 83551670
$ws.Range("I94").Value = 250025000
$ws.Range("K94").Value = 250025000
$ws.Range("M94").Value = -250024549
$ws.Range("H100").Value = 8331382
$ws.Range("I100").Value = 10006055
$ws.Range("K100").Value = 10006055
$ws.Range("M100").Value = -10005514
$ws.Range("H112").Value = 36412.855
$ws.Range("J112").Value = 37395.59
$ws.Range("L112").Value = 112186.77
$ws.Range("N112").Value = -114402.77
$ws.Range("H113").Value = 11875.667
$ws.Range("I113").Value = 14433.154
$ws.Range("K113").Value = 14433.154
$ws.Range("M113").Value = -11179.154
$ws.Range("H115").Value = 370
$ws.Range("I115").Value = 370
$ws.Range("K115").Value = 1110
$ws.Range("M115").Value = 457
$ws.Range("H132").Value = 7796.8125
$ws.Range("I132").Value = 7835
$ws.Range("K132").Value = 23505
$ws.Range("M132").Value = -20975

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 2243.3333
$ws.Range("I61").Value = 2243.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2243.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2031.3333
$ws.Range("H88").Value = 111113370
$ws.Range("I88").Value = 2000.5
$ws.Range("K88").Value = 2000.5
$ws.Range("M88").Value = -1594.5
$ws.Range("H91").Value = 111113370
$ws.Range("I91").Value = 2000.5
$ws.Range("K91").Value = 2000.5
$ws.Range("M91").Value = -596.5
$ws.Range("H97").Value = 22260652
$ws.Range("I97").Value = 35054
$ws.Range("J97").Value = 33373452
$ws.Range("K97").Value = 35054
$ws.Range("L97").Value = 33373452
$ws.Range("M97").Value = -34558
$ws.Range("N97").Value = -33374444
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 0
$ws.Range("H122").Value = 564129.75
$ws.Range("I122").Value = 3203
$ws.Range("J122").Value = 1405519.9
$ws.Range("K122").Value = 9609
$ws.Range("L122").Value = 4216559.699999999
$ws.Range("M122").Value = -7159
$ws.Range("N122").Value = -4221459.699999999
$ws.Range("H132").Value = 2465.1282
$ws.Range("I132").Value = 2133.1333
$ws.Range("K132").Value = 6399.3999
$ws.Range("M132").Value = -3869.3999
$ws.Range("H136").Value = 2243.3333
$ws.Range("I136").Value = 2243.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6729.999899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -4179.999899999999

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 25704.092
$ws.Range("I94").Value = 44181.332
$ws.Range("J94").Value = 3531.4
$ws.Range("K94").Value = 44181.332
$ws.Range("L94").Value = 3531.4
$ws.Range("M94").Value = -43730.332
$ws.Range("N94").Value = -4433.4
$ws.Range("H123").Value = 77000
$ws.Range("J123").Value = 77000
$ws.Range("L123").Value = 77000
$ws.Range("N123").Value = -86800

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H25").Value = 4714.2856
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 4923.077
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 4923.077
$ws.Range("M25").Value = -1826
$ws.Range("N25").Value = -5271.077
$ws.Range("H43").Value = 76851
$ws.Range("J43").Value = 76851
$ws.Range("L43").Value = 76851
$ws.Range("N43").Value = -77219
$ws.Range("H101").Value = 76851
$ws.Range("J101").Value = 76851
$ws.Range("L101").Value = 76851
$ws.Range("N101").Value = -83341
$ws.Range("H105").Value = 101648.664
$ws.Range("I105").Value = 132631.38
$ws.Range("K105").Value = 132631.38
$ws.Range("M105").Value = -130884.38

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 30731454
$ws.Range("I4").Value = 31635262
$ws.Range("K4").Value = 94905786
$ws.Range("M4").Value = -94905674
$ws.Range("H68").Value = 16135938
$ws.Range("J68").Value = 21747552
$ws.Range("L68").Value = 65242656
$ws.Range("N68").Value = -65244278
$ws.Range("H71").Value = 16135938
$ws.Range("J71").Value = 21747552
$ws.Range("L71").Value = 195727968
$ws.Range("N71").Value = -195736080
$ws.Range("H74").Value = 100
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 100
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 20899324
$ws.Range("I132").Value = 2975
$ws.Range("J132").Value = 27864774
$ws.Range("K132").Value = 26775
$ws.Range("L132").Value = 250782966
$ws.Range("M132").Value = -24245
$ws.Range("N132").Value = -250788026

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H97").Value = 11640.7
$ws.Range("I97").Value = 14939.6
$ws.Range("J97").Value = 1744
$ws.Range("K97").Value = 14939.6
$ws.Range("L97").Value = 1744
$ws.Range("M97").Value = -14443.6
$ws.Range("N97").Value = -2736
$ws.Range("H102").Value = 7545.269
$ws.Range("I102").Value = 7798.64
$ws.Range("K102").Value = 7798.64
$ws.Range("M102").Value = -6176.64
$ws.Range("H107").Value = 622.5
$ws.Range("I107").Value = 563.3333
$ws.Range("K107").Value = 563.3333
$ws.Range("M107").Value = 1356.6667
$ws.Range("H132").Value = 4353.4863
$ws.Range("I132").Value = 3436.7
$ws.Range("K132").Value = 10310.1
$ws.Range("M132").Value = -7780.099999999999

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H93").Value = 7555.2856
$ws.Range("I93").Value = 7555.2856
$ws.Range("K93").Value = 7555.2856
$ws.Range("M93").Value = -6307.2856

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H41").Value = 21833.334
$ws.Range("J41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30780
$ws.Range("H81").Value = 36149.832
$ws.Range("I81").Value = 51225
$ws.Range("J81").Value = 5999.5
$ws.Range("K81").Value = 102450
$ws.Range("L81").Value = 11999
$ws.Range("M81").Value = -101389
$ws.Range("N81").Value = -14121
$ws.Range("H84").Value = 36149.832
$ws.Range("I84").Value = 51225
$ws.Range("J84").Value = 5999.5
$ws.Range("K84").Value = 512250
$ws.Range("L84").Value = 59995
$ws.Range("M84").Value = -506946
$ws.Range("N84").Value = -70603
$ws.Range("H100").Value = 21766.965
$ws.Range("I100").Value = 13286.044
$ws.Range("J100").Value = 60779.2
$ws.Range("K100").Value = 26572.088
$ws.Range("L100").Value = 121558.4
$ws.Range("M100").Value = -26031.088
$ws.Range("N100").Value = -122640.4
$ws.Range("H136").Value = 3100.28
$ws.Range("I136").Value = 2358.9375
$ws.Range("K136").Value = 7076.8125
$ws.Range("M136").Value = -4526.8125
